$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$NewValue
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

Set-TextValue "D2" "61.695.02"
Set-TextValue "E2" "  -1.54%  "
Set-TextValue "D3" "3.046.46"
Set-TextValue "E3" "  -4.30%  "
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "583.03"
Set-TextValue "E5" "  -1.05%  "
Set-TextValue "D6" "130.54"
Set-TextValue "E6" "  -4.00%  "
Set-TextValue "E7" "  +0.03%  "
Set-TextValue "D8" "3.047.92"
Set-TextValue "E8" "  -4.17%  "
Set-TextValue "D9" "0.504"
Set-TextValue "E9" "  -0.54%  "
Set-TextValue "E10" "  -2.22%  "
Set-TextValue "D11" "5.29"
Set-TextValue "E11" "  +0.17%  "
Set-TextValue "E12" "  -2.93%  "
Set-TextValue "D13" "0.0000233"
Set-TextValue "E13" "  -1.07%  "
Set-TextValue "D14" "33.65"
Set-TextValue "E14" "  +0.70%  "
Set-TextValue "E15" "  +0.90%  "
Set-TextValue "D16" "3.553.46"
Set-TextValue "E16" "  -4.09%  "
Set-TextValue "D17" "61.747.85"
Set-TextValue "E17" "  -1.42%  "
Set-TextValue "D18" "3.048.65"
Set-TextValue "E18" "  -4.05%  "
Set-TextValue "D19" "6.37"
Set-TextValue "E19" "  -2.72%  "
Set-TextValue "D20" "449.41"
Set-TextValue "E20" "  -1.61%  "
Set-TextValue "D21" "13.53"
Set-TextValue "E21" "  -2.86%  "
Set-TextValue "D22" "0.672"
Set-TextValue "E22" "  -4.54%  "
Set-TextValue "D23" "7.35"
Set-TextValue "E23" "  -3.85%  "
Set-TextValue "D24" "81.10"
Set-TextValue "E24" "  -2.86%  "
Set-TextValue "D25" "12.84"
Set-TextValue "E25" "  -4.26%  "
Set-TextValue "E26" "  +0.07%  "
Set-TextValue "E27" "  -0.15%  "
Set-TextValue "D28" "2.57"
Set-TextValue "E28" "  -4.76%  "
Set-TextValue "E29" "  -0.65%  "
Set-TextValue "D30" "7.41"
Set-TextValue "E30" "  -4.73%  "
Set-TextValue "D31" "6.44"
Set-TextValue "E31" "  -6.50%  "
Set-TextValue "D32" "25.92"
Set-TextValue "E32" "  -5.44%  "
Set-TextValue "D33" "0.0974"
Set-TextValue "E33" "  -6.39%  "
Set-TextValue "D34" "2.33"
Set-TextValue "E34" "  -2.38%  "
Set-TextValue "E35" "  -5.85%  "
Set-TextValue "E36" "  -3.21%  "
Set-TextValue "D37" "50.39"
Set-TextValue "E37" "  -1.29%  "
Set-TextValue "D38" "0.0₃0695"
Set-TextValue "E38" "  -0.82%  "
Set-TextValue "D39" "0.0375"
Set-TextValue "E39" "  -2.90%  "
Set-TextValue "D40" "7.93"
Set-TextValue "E40" "  -1.06%  "
Set-TextValue "B41" "Kaspa"
Set-TextValue "C41" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.109"
Set-TextValue "E41" "  -3.03%  "
Set-TextValue "B42" "Bittensor"
Set-TextValue "C42" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D42" "381.27"
Set-TextValue "E42" "  -6.92%  "
Set-TextValue "E43" "  -7.38%  "
Set-TextValue "D44" "2.702.05"
Set-TextValue "E44" "  -5.32%  "
Set-TextValue "E45" "  +0.00%  "
Set-TextValue "D46" "123.87"
Set-TextValue "E46" "  -0.99%  "
Set-TextValue "D47" "0.240"
Set-TextValue "E47" "  -3.86%  "
Set-TextValue "B48" "Arweave"
Set-TextValue "C48" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D48" "34.18"
Set-TextValue "E48" "  -6.32%  "
Set-TextValue "B49" "Fetch.AI"
Set-TextValue "C49" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D49" "2.02"
Set-TextValue "E49" "  -5.60%  "
Set-TextValue "E50" "  -2.30%  "
Set-TextValue "D51" "23.95"
Set-TextValue "E51" "  -6.17%  "
